# Updates the "10/4/2023" date placeholders (slide master, all 11 slide
# layouts, and the notes master) to "10/14/2024", and updates the
# copyright year / re-splits the credit line on slide 1 from
# "Slides  ©2023 Rose bohrer, used for cs 4536/536 at ..." to
# "Slides  ©2024 Rose bohrer, used for cs 4536/536 at ...".

$p = $ppt.ActivePresentation

$oldDate = "10/4/2023"
$newDate = "10/14/2024"

# --- Slide master date placeholder -----------------------------------
$masterDateShape = $p.SlideMaster.Shapes.Item("Date Placeholder 3")
$masterDateShape.TextFrame.TextRange.Text = $newDate

# --- All 11 slide layouts' date placeholders --------------------------
$layoutDateShapeNames = @(
    "Date Placeholder 3",  # 1  Title Slide
    "Date Placeholder 3",  # 2  Title and Content
    "Date Placeholder 3",  # 3  Section Header
    "Date Placeholder 4",  # 4  Two Content
    "Date Placeholder 6",  # 5  Comparison
    "Date Placeholder 2",  # 6  Title Only
    "Date Placeholder 6",  # 7  Blank
    "Date Placeholder 4",  # 8  Content with Caption
    "Date Placeholder 4",  # 9  Picture with Caption
    "Date Placeholder 3",  # 10 Title and Vertical Text
    "Date Placeholder 3"   # 11 Vertical Title and Text
)

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    $dateShape = $layout.Shapes.Item($layoutDateShapeNames[$i - 1])
    $dateShape.TextFrame.TextRange.Text = $newDate
}

# --- Notes master date placeholder ------------------------------------
$notesDateShape = $p.NotesMaster.Shapes.Item("Date Placeholder 2")
$notesDateShape.TextFrame.TextRange.Text = $newDate

# --- Slide 1 subtitle: year bump + credit-line retouch -----------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item("Subtitle 2")
$tr = $subtitle.TextFrame.TextRange

# "Slides  ©2023 Rose " -> split into "Slides  " / "©2024 " / "Rose "
$yearRange = $tr.Characters(9, 6)
$yearRange.Text = "©2024 "

# "bohrer" + ", used for cs 4536/536 at " -> "bohrer," + " used for cs 4536/536 at "
# (moves the comma so it is attached to "bohrer" instead of starting the next run)
$commaRange = $tr.Characters(20, 7)
$commaRange.Text = $commaRange.Text
